$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BidAISimulation")

# Update input values on the BidAISimulation sheet
$ws.Range("I3").Value = 4
$ws.Range("I4").Value = 4
$ws.Range("I6").Value = 0
$ws.Range("E7").Value = 1
$ws.Range("I7").Value = 1

# Recalculate dependent formulas (I8, I12, H13, I13, etc.)
$excel.CalculateFull()

# Update the active cell selection to match the saved view state
$ws.Activate()
$ws.Range("T9").Select()
